# Daily attendance processing
# Normalizes the "Recorded By" column (G) so that "System" always appears
# first in the comma-separated list of recorders, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Entries that already start with "System", or that don't mention it at
# all, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$suffix = ", System"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        if ($val.EndsWith($suffix) -and -not $val.StartsWith("System")) {
            $rest = $val.Substring(0, $val.Length - $suffix.Length)
            $cell.Value = "System, " + $rest
        }
    }
}
